# Apply updated crypto price/volume data to worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.255.85"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.56%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.275.18"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.74%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "546.56"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.70%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.22"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.60%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.12%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.57%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.45"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.59%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.117"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.99%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.435"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.92%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.841.27"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.85%  "

# Row 14
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000178"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.59%  "

# Row 15
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.61"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.79%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.335.71"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.69%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.275.64"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.32%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.36"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.55%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.62"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.83%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.45"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.90%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "379.58"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.32%  "

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.27%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.532"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.05%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.26"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.15%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.173"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.39%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.82"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.71%  "

# Row 27
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0966"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +7.68%  "

# Row 28
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.41%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.95"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.02%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.77"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.96%  "

# Row 31
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.29"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +8.53%  "

# Row 32
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "RenderToken"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.23"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.21%  "

# Row 33
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.44"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.22%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.70"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.62%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "158.91"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.92%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.46"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +6.99%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.58"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.80%  "

# Row 38
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.80"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +6.57%  "

# Row 39
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.808.73"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.13%  "

# Row 40
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0723"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.85%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0318"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +7.73%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.29"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.14%  "

# Row 43
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.739"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.63%  "

# Row 44
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.08"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.74%  "

# Row 45
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "ONDO"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.03"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.90%  "

# Row 46
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.105"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.54%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.79"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +5.28%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.321.27"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.86%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.26"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.53%  "

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +4.62%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "277.23"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +6.37%  "
